$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "b79f5a86"
$ws.Range("B2").Value = "b2ea9becb2e0"
$ws.Range("C2").Value = "Smith"
$ws.Range("D2").Value = "Main Clinic"
$ws.Range("E2").Value = "2025-09-06T09:00:00+00:00"
$ws.Range("F2").Value = 60
